$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Update the "Produzione" value for row 24 (CANNARELLA GIUSEPPE) from 35.99 to 35.
$ws.Range("C24").Value = 35

# Re-enter the "Avanzamento €/h" formula across E2:E25 in one shot so Excel
# consolidates it into a shared formula group (master at E2, si="0"),
# matching the recalculated dependent value in E24.
$ws.Range("E2:E25").Formula = "=C2-(C2*D2)/100"

# Move the active selection to C26, as recorded in the saved view state.
$ws.Range("C26").Select() | Out-Null
